# The source ads-scraping sheet picked up a second listing and a couple of
# numeric-looking fields on the first row now parse as real numbers instead
# of text - mirror that here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: a handful of cells that used to be scraped as text now come in
# as genuine numbers.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 30598523    # A2 Ad ID
$ws.Cells.Item(2, 4).Value = 29          # D2 Number
$ws.Cells.Item(2, 6).Value = 250         # F2 Agency Price
$ws.Cells.Item(2, 12).Value = 2          # L2 ოთახი
$ws.Cells.Item(2, 13).Value = 1          # M2 საძინებელი
$ws.Cells.Item(2, 14).Value = 1          # N2 სართული
$ws.Cells.Item(2, 15).Value = 7          # O2 სართულიანობა
$ws.Cells.Item(2, 16).Value = 1          # P2 სველი წერტილი

# J2 (Comment) now exists but is simply blank text - enter it the way
# Excel itself records an intentionally-empty text cell (leading
# apostrophe forces text, then drop the formatting it implies).
$ws.Cells.Item(2, 10).Value = "'"
$ws.Cells.Item(2, 10).ClearFormats()

# ---------------------------------------------------------------------
# Row 3: a brand-new scraped listing.
# ---------------------------------------------------------------------
# A3 looks like an integer ("31024583") but the source keeps it as text -
# a bare numeric-looking string auto-coerces to a number via .Value, so
# force text with a leading apostrophe (then drop the quote-prefix format
# it implies) just like J2 above.
$ws.Cells.Item(3, 1).Value = "'31024583"
$ws.Cells.Item(3, 1).ClearFormats()

$ws.Cells.Item(3, 2).Value = "ქირავდება 2 ოთახიანი ბინა ისანში"
$ws.Cells.Item(3, 3).Value = "ბერი გაბრიელ სალოსის გამზ."

$ws.Cells.Item(3, 4).Value = "'"
$ws.Cells.Item(3, 4).ClearFormats()

$ws.Cells.Item(3, 5).Value = "1,200 ₾"

# F3 ("1100") likewise needs to stay text.
$ws.Cells.Item(3, 6).Value = "'1100"
$ws.Cells.Item(3, 6).ClearFormats()

$ws.Cells.Item(3, 7).Value = "558 46 27 74"
$ws.Cells.Item(3, 8).Value = "nini"
$ws.Cells.Item(3, 9).Value = "ს ა ს წ რ ა ფ ო დ !!! ისანში ბერი გაბრიელ სალოსის გამზირზე ქირავდება 2 ოთახიანი ბინა 1 საძინებელი ოთახით,ავეჯით,ტექნიკით,ცენტრალური გათბობის სიტემით,კეთილმოწყობილი შლაგბაუმიანი ეზოთი,კოდირებული კარით,ვიდეომონიტორინგით სადარბაზოში,მიმდებარედ ყველანაირი ობიექტით."
$ws.Cells.Item(3, 10).Value = "sdadasd"
$ws.Cells.Item(3, 11).Value = "56 მ²"

# L3..P3 ("2","1","10","12","1") are the same numeric-looking-but-text
# case as A3/F3 above.
$ws.Cells.Item(3, 12).Value = "'2"
$ws.Cells.Item(3, 12).ClearFormats()
$ws.Cells.Item(3, 13).Value = "'1"
$ws.Cells.Item(3, 13).ClearFormats()
$ws.Cells.Item(3, 14).Value = "'10"
$ws.Cells.Item(3, 14).ClearFormats()
$ws.Cells.Item(3, 15).Value = "'12"
$ws.Cells.Item(3, 15).ClearFormats()
$ws.Cells.Item(3, 16).Value = "'1"
$ws.Cells.Item(3, 16).ClearFormats()

$ws.Cells.Item(3, 17).Value = "ახალი რემონტით"
$ws.Cells.Item(3, 18).Value = "ახალი აშენებული"
$ws.Cells.Item(3, 19).Value = "უძრავი ქონება"
$ws.Cells.Item(3, 20).Value = "ბინა"
$ws.Cells.Item(3, 21).Value = "ქირავდება"
$ws.Cells.Item(3, 22).Value = "კი"
$ws.Cells.Item(3, 23).Value = "კი"
$ws.Cells.Item(3, 24).Value = "არა"
$ws.Cells.Item(3, 25).Value = "კი"
$ws.Cells.Item(3, 26).Value = "კი"
$ws.Cells.Item(3, 27).Value = "კი"
$ws.Cells.Item(3, 28).Value = "კი"
$ws.Cells.Item(3, 29).Value = "არა"
$ws.Cells.Item(3, 30).Value = "კი"
$ws.Cells.Item(3, 31).Value = "კი"
$ws.Cells.Item(3, 32).Value = "კი"
$ws.Cells.Item(3, 33).Value = "კი"
$ws.Cells.Item(3, 34).Value = "კი"
$ws.Cells.Item(3, 35).Value = "კი"
$ws.Cells.Item(3, 36).Value = "არა"
$ws.Cells.Item(3, 37).Value = "არა"
$ws.Cells.Item(3, 38).Value = "არა"
$ws.Cells.Item(3, 39).Value = "კი"
$ws.Cells.Item(3, 40).Value = "კი"

# AO3 (Final URL) - this listing's final URL hasn't been resolved yet,
# so it's present but blank, same trick as J2/D3 above.
$ws.Cells.Item(3, 41).Value = "'"
$ws.Cells.Item(3, 41).ClearFormats()
